$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.803.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.811.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.33"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4321"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3709"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07240"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8671"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.98%  "

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.81"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.56%  "

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.951.08"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +6.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.688"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.357"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06913"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.010"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "80.47"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008918"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.79%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.27"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.852.54"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.215"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.16"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.219.83"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +8.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.65"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.874"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.84%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.201"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.907"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +15.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.32"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08937"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7551"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.165"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.435"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.807"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.009"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.126"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05222"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01926"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5076"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1650"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.665"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.566"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +10.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.294"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.48%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.44"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.20%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.38"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.005"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.654"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.24%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06280"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.51%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4555"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.801"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.66%  "
